# "Generate Report for Handback"
#
# For each localized-language sheet (zh-cn, de-de):
#   - Status (col C) moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" (this is a shared-string value,
#     so simply overwriting either cell's text updates the single shared
#     string and every cell that referenced it).
#   - Latest Target File (col F) and Latest Handback File (col G) get
#     populated with hyperlinked filenames (the handed-off source .md and
#     the translated .xlf respectively) for each of the two data rows.
#   - Latest Handback DateTime (col H) moves from the zero date to the
#     actual handback timestamp for that language.

$wb = $excel.ActiveWorkbook

# ---- Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
# Touching it once on the Overview sheet updates the shared string that
# every other reference (zh-cn!C2/C3, de-de!C2/C3, Overview!B/C) points to.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"

$languages = @(
    @{
        SheetName = "zh-cn"
        HandbackDateTime = "2016-03-18 04:22:41"
        Row2 = @{
            TargetFile  = "016ed19e-e6b8-4d47-b037-60934746d37f.md"
            TargetUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/ed3713c25a602242e3eeb9d6fd1754aa014467ea/e2e/016ed19e-e6b8-4d47-b037-60934746d37f.md"
            HandbackFile = "016ed19e-e6b8-4d47-b037-60934746d37f.e558124005ca1c2f6a3906d810bb2e52dbecfa95.zh-cn.xlf"
            HandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ed82bf863ed7037380c7d589da17b3331660ca8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/016ed19e-e6b8-4d47-b037-60934746d37f.e558124005ca1c2f6a3906d810bb2e52dbecfa95.zh-cn.xlf"
        }
        Row3 = @{
            TargetFile  = "9bfdd7b6-0537-41ef-9a45-339b14cae5cd.md"
            TargetUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/ed3713c25a602242e3eeb9d6fd1754aa014467ea/e2e/9bfdd7b6-0537-41ef-9a45-339b14cae5cd.md"
            HandbackFile = "9bfdd7b6-0537-41ef-9a45-339b14cae5cd.295ebee6c9499d6623ae6341a32e1d47f94e548a.zh-cn.xlf"
            HandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0ed82bf863ed7037380c7d589da17b3331660ca8/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9bfdd7b6-0537-41ef-9a45-339b14cae5cd.295ebee6c9499d6623ae6341a32e1d47f94e548a.zh-cn.xlf"
        }
    },
    @{
        SheetName = "de-de"
        HandbackDateTime = "2016-03-18 04:22:46"
        Row2 = @{
            TargetFile  = "016ed19e-e6b8-4d47-b037-60934746d37f.md"
            TargetUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/ed3713c25a602242e3eeb9d6fd1754aa014467ea/e2e/016ed19e-e6b8-4d47-b037-60934746d37f.md"
            HandbackFile = "016ed19e-e6b8-4d47-b037-60934746d37f.e558124005ca1c2f6a3906d810bb2e52dbecfa95.de-de.xlf"
            HandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cacd313411ae6009b84b4588b68eb932172b9927/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/016ed19e-e6b8-4d47-b037-60934746d37f.e558124005ca1c2f6a3906d810bb2e52dbecfa95.de-de.xlf"
        }
        Row3 = @{
            TargetFile  = "9bfdd7b6-0537-41ef-9a45-339b14cae5cd.md"
            TargetUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/ed3713c25a602242e3eeb9d6fd1754aa014467ea/e2e/9bfdd7b6-0537-41ef-9a45-339b14cae5cd.md"
            HandbackFile = "9bfdd7b6-0537-41ef-9a45-339b14cae5cd.295ebee6c9499d6623ae6341a32e1d47f94e548a.de-de.xlf"
            HandbackUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cacd313411ae6009b84b4588b68eb932172b9927/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9bfdd7b6-0537-41ef-9a45-339b14cae5cd.295ebee6c9499d6623ae6341a32e1d47f94e548a.de-de.xlf"
        }
    }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.SheetName)

    # Status column (just re-asserting the post-handback text; harmless if
    # it was already updated via the shared string above).
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
    $ws.Range("C3").Value = "Handed back: in sync with en-US"

    # Row 2 (016ed19e... file)
    $ws.Range("F2").Value = $lang.Row2.TargetFile
    $ws.Hyperlinks.Add($ws.Range("F2"), $lang.Row2.TargetUrl, "", "", $lang.Row2.TargetFile) | Out-Null

    $ws.Range("G2").Value = $lang.Row2.HandbackFile
    $ws.Hyperlinks.Add($ws.Range("G2"), $lang.Row2.HandbackUrl, "", "", $lang.Row2.HandbackFile) | Out-Null

    $ws.Range("H2").Value = $lang.HandbackDateTime

    # Row 3 (9bfdd7b6... file)
    $ws.Range("F3").Value = $lang.Row3.TargetFile
    $ws.Hyperlinks.Add($ws.Range("F3"), $lang.Row3.TargetUrl, "", "", $lang.Row3.TargetFile) | Out-Null

    $ws.Range("G3").Value = $lang.Row3.HandbackFile
    $ws.Hyperlinks.Add($ws.Range("G3"), $lang.Row3.HandbackUrl, "", "", $lang.Row3.HandbackFile) | Out-Null

    $ws.Range("H3").Value = $lang.HandbackDateTime
}
